$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-09 -> 2023-09-10, i.e. 45178 -> 45179) for every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 339) { $lastRow = 339 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cur = $ws.Cells.Item($row, 3).Value2
    if ($cur -eq 45178) {
        $ws.Cells.Item($row, 3).Value = 45179
    }
}
